$d = $word.ActiveDocument

function Get-CellText($cell) {
    return $cell.Range.Text.TrimEnd([char]7, [char]13, [char]12)
}

# ------------------------------------------------------------------
# 1) Remove the stray _GoBack bookmark that currently sits at the end
#    of the intro paragraph (after the NFS description text).
# ------------------------------------------------------------------
try {
    $goBack = $d.Bookmarks("_GoBack")
    $goBack.Delete()
} catch {
    # no pre-existing _GoBack bookmark; nothing to remove
}

# ------------------------------------------------------------------
# 2) Add the missing bullet text under "Inclusions" (numId 6) and
#    "Exclusions" (numId 4) - the first empty ListParagraph after each
#    heading's intro sentence.
# ------------------------------------------------------------------
$inclusionsDone = $false
$exclusionsDone = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $txt = $p.Range.Text.TrimEnd([char]7, [char]13, [char]12)

    if (-not $inclusionsDone -and $txt -eq "The following items will be included in the scope of work:") {
        $next = $d.Paragraphs.Item($i + 1)
        if ($next.Range.Text.TrimEnd([char]7, [char]13, [char]12) -eq "") {
            $next.Range.InsertAfter("Messages sent from client to server.")
        }
        $inclusionsDone = $true
    }

    if (-not $exclusionsDone -and $txt -eq "The following items will not be included in the scope of work:") {
        $next = $d.Paragraphs.Item($i + 1)
        if ($next.Range.Text.TrimEnd([char]7, [char]13, [char]12) -eq "") {
            $next.Range.InsertAfter("Server response messages")
        }
        $exclusionsDone = $true
    }
}

# ------------------------------------------------------------------
# 3) Schedule table: "Build" work item duration 3 weeks -> 2 weeks
# ------------------------------------------------------------------
$scheduleTbl = $d.Tables.Item(1)
for ($r = 1; $r -le $scheduleTbl.Rows.Count; $r++) {
    $whatCell = $scheduleTbl.Cell($r, 3)
    if ((Get-CellText $whatCell) -eq "Build") {
        $whenCell = $scheduleTbl.Cell($r, 4)
        $whenStart = $whenCell.Range.Start
        $digit = $d.Range($whenStart, $whenStart + 1)
        if ($digit.Text -eq "3") {
            $digit.Text = "2"
        }
    }
}

# ------------------------------------------------------------------
# 4) Hours Breakdown table: "Build Pit" 15 days -> 10 days, and
#    "TOTAL HOURS" 22 days -> 17 days (re-homing the _GoBack bookmark
#    to sit right after the new total number).
# ------------------------------------------------------------------
$hoursTbl = $d.Tables.Item(2)
for ($r = 1; $r -le $hoursTbl.Rows.Count; $r++) {
    $labelCell = $hoursTbl.Cell($r, 1)
    $label = Get-CellText $labelCell

    if ($label -eq "Build Pit") {
        $hoursCell = $hoursTbl.Cell($r, 2)
        $hoursStart = $hoursCell.Range.Start
        $digit = $d.Range($hoursStart + 1, $hoursStart + 2)
        if ($digit.Text -eq "5") {
            $digit.Text = "0"
        }
    }

    if ($label -eq "TOTAL HOURS") {
        $hoursCell = $hoursTbl.Cell($r, 2)
        $hoursStart = $hoursCell.Range.Start
        $numRng = $d.Range($hoursStart, $hoursStart + 2)
        if ($numRng.Text -eq "22") {
            $numRng.Text = "17"
        }
        $newBmRange = $d.Range($hoursStart + 2, $hoursStart + 2)
        $d.Bookmarks.Add("_GoBack", $newBmRange)
    }
}
